$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.024.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "'1.645.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'215.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "'0.5226"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.06360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "'20.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").Value = "'0.07675"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "'1.645.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'4.421"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'1.867.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "'0.5537"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "'0.0₅8320"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").Value = "'64.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "'26.024.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D20").Value = "'4.717"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "'188.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "'10.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'6.258"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'144.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("D26").Value = "'0.1220"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'7.403"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "'15.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'1.388"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "'0.05959"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.63%  "
$ws.Range("D31").Value = "'1.264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").Value = "'3.397"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'3.402"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").Value = "'1.652"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'0.9954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'2.392"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").Value = "'2.753"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "'0.5632"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.19%  "
$ws.Range("D39").Value = "'0.01607"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'5.852"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").Value = "'0.8548"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'1.026.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.76%  "
$ws.Range("D44").Value = "'98.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'1.794.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "'55.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D49").Value = "'8.048"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'0.05145"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "'0.4215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.53%  "
